$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 42 ("miss"), pushing
# "miss" and "END" down to rows 44 and 45. This is the new "self"/"other"
# forced-targeting keyword pair added to the style-modifier block.
$ws.Rows("42:43").Insert()

# Populate the "other" row (43) first, then the "self" row (42), so the
# shared-string table records "other"/"forces spell to target other"
# ahead of "self"/"forces spell to target self" (matching authoring order).
$ws.Range("A43").Value = "other"
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = "forces spell to target other"
$ws.Range("D43").Value = "anim_spell_empower"
$ws.Range("E43").Value = "sfx_mystery_effect"
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 1
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 1
$ws.Range("M43").Value = 0
$ws.Range("N43").Value = 1
$ws.Range("O43").Value = "MT"

$ws.Range("A42").Value = "self"
$ws.Range("B42").Value = 0
$ws.Range("C42").Value = "forces spell to target self"
$ws.Range("D42").Value = "anim_spell_empower"
$ws.Range("E42").Value = "sfx_mystery_effect"
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 1
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 1
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 1
$ws.Range("M42").Value = 0
$ws.Range("N42").Value = 1
$ws.Range("O42").Value = "mS"

# Re-assigning .Value above reset the quote-prefix formatting that
# Rows.Insert had copied down into row 42/43's C and F cells; restore it
# by copying the format from the row above (still intact).
$ws.Range("C41").Copy()
$ws.Range("C42:C43").PasteSpecial(-4122)
$ws.Range("F41").Copy()
$ws.Range("F42:F43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the editing session's final selection (landed on H43), like
# the source workbook did.
[void]$ws.Range("H43").Select()
